$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Q8" header label in J1
$ws.Range("J1").Value = "Q8"

# Update data values in B2:J16 to reflect recomputed naive error table
$ws.Range("B2").Value = 0.5556414999948345
$ws.Range("C2").Value = -0.3970860714793787
$ws.Range("D2").Value = -0.6169817575957095
$ws.Range("E2").Value = 1.424879646726211
$ws.Range("F2").Value = 1.650808050689989
$ws.Range("G2").Value = 0.3735921536389384
$ws.Range("H2").Value = 0.6866231187695604
$ws.Range("B3").Value = 0.7219671938530607
$ws.Range("C3").Value = 0.5020715077367299
$ws.Range("D3").Value = 2.54393291205865
$ws.Range("E3").Value = 2.769861316022428
$ws.Range("F3").Value = 1.492645418971378
$ws.Range("G3").Value = 1.805676384102
$ws.Range("B4").Value = -1.186928492263277
$ws.Range("C4").Value = 0.8549329120586431
$ws.Range("D4").Value = 1.080861316022421
$ws.Range("E4").Value = -0.1963545810286291
$ws.Range("F4").Value = 0.1166763841019929
$ws.Range("G4").Value = -0.3857387297557864
$ws.Range("H4").Value = 0.2658945176530381
$ws.Range("I4").Value = 0.0729948986110571
$ws.Range("J4").Value = -0.4377465499738378
$ws.Range("B5").Value = 0.5039362906915836
$ws.Range("C5").Value = 0.7298646946553617
$ws.Range("D5").Value = -0.5473512023956886
$ws.Range("E5").Value = -0.2343202372650666
$ws.Range("F5").Value = -0.7367353511228459
$ws.Range("G5").Value = -0.08510210371402144
$ws.Range("H5").Value = -0.2780017227560024
$ws.Range("I5").Value = -0.7887431713408973
$ws.Range("B6").Value = 0.4898613160224272
$ws.Range("C6").Value = -0.7873545810286231
$ws.Range("D6").Value = -0.4743236158980011
$ws.Range("E6").Value = -0.9767387297557804
$ws.Range("F6").Value = -0.3251054823469559
$ws.Range("G6").Value = -0.5180051013889369
$ws.Range("H6").Value = -1.028746549973832
$ws.Range("B7").Value = -0.298354581028633
$ws.Range("C7").Value = 0.01467638410198902
$ws.Range("D7").Value = -0.4877387297557902
$ws.Range("E7").Value = 0.1638945176530342
$ws.Range("F7").Value = -0.02900510138894677
$ws.Range("G7").Value = -0.5397465499738416
$ws.Range("B8").Value = 0.4656763841019966
$ws.Range("C8").Value = -0.03673872975578271
$ws.Range("D8").Value = 0.6148945176530418
$ws.Range("E8").Value = 0.4219948986110608
$ws.Range("F8").Value = -0.08874654997383413
$ws.Range("G8").Value = 0.116654182623364
$ws.Range("H8").Value = -0.2898317593399469
$ws.Range("I8").Value = 0.1452723979283945
$ws.Range("B9").Value = -0.116738729755781
$ws.Range("C9").Value = 0.5348945176530435
$ws.Range("D9").Value = 0.3419948986110625
$ws.Range("E9").Value = -0.1687465499738324
$ws.Range("F9").Value = 0.03665418262336569
$ws.Range("G9").Value = -0.3698317593399452
$ws.Range("H9").Value = 0.06527239792839619
$ws.Range("B10").Value = 0.5148945176530333
$ws.Range("C10").Value = 0.3219948986110523
$ws.Range("D10").Value = -0.1887465499738426
$ws.Range("E10").Value = 0.0166541826233555
$ws.Range("F10").Value = -0.3898317593399554
$ws.Range("G10").Value = 0.045272397928386
$ws.Range("B11").Value = 0.221994898611058
$ws.Range("C11").Value = -0.2887465499738369
$ws.Range("D11").Value = -0.08334581737663882
$ws.Range("E11").Value = -0.4898317593399497
$ws.Range("F11").Value = -0.05472760207160832
$ws.Range("B12").Value = -0.1887465499738426
$ws.Range("C12").Value = 0.0166541826233555
$ws.Range("D12").Value = -0.3898317593399554
$ws.Range("E12").Value = 0.045272397928386
$ws.Range("B13").Value = 0.2366557702529377
$ws.Range("C13").Value = -0.1698301717103732
$ws.Range("D13").Value = 0.2652739855579682
$ws.Range("B14").Value = -0.3898340444052479
$ws.Range("C14").Value = 0.04527011286309346
$ws.Range("B15").Value = 0.1452758398526868
